$d = $word.ActiveDocument

# Merge the split "<id>p035v_1</id>" runs into a single run with the
# combined text (keeping the formatting of the first run, matching how
# Word's Find & Replace across multiple runs collapses into one run).
$d.Content.Find.Execute("<id>p035v_1</id>", $false, $false, $false, $false, $false, `
                         $true, 1, $false, "<id>p035v_1</id>", 2)

# Same for the second occurrence "<id>p035v_2</id>".
$d.Content.Find.Execute("<id>p035v_2</id>", $false, $false, $false, $false, $false, `
                         $true, 1, $false, "<id>p035v_2</id>", 2)
